# Fruta / hortaliza, semanal
# Weekly data refresh: reassigns the Fecha (D) and Volumen/Precio
# (J, K, L, M, P) columns for rows 2-31 on the active sheet to the
# new reported values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 45041
$ws.Cells.Item(2, 10).Value = 1160
$ws.Cells.Item(2, 11).Value = 2000
$ws.Cells.Item(2, 12).Value = 2500
$ws.Cells.Item(2, 13).Value = 2250
$ws.Cells.Item(2, 16).Value = 750

$ws.Cells.Item(3, 4).Value = 44827
$ws.Cells.Item(3, 10).Value = 1200

$ws.Cells.Item(4, 4).Value = 45013

$ws.Cells.Item(5, 4).Value = 45034
$ws.Cells.Item(5, 10).Value = 1100
$ws.Cells.Item(5, 11).Value = 2000
$ws.Cells.Item(5, 12).Value = 2500
$ws.Cells.Item(5, 13).Value = 2250
$ws.Cells.Item(5, 16).Value = 750

$ws.Cells.Item(6, 4).Value = 44951
$ws.Cells.Item(6, 10).Value = 800

$ws.Cells.Item(7, 4).Value = 45020
$ws.Cells.Item(7, 10).Value = 1200

$ws.Cells.Item(9, 4).Value = 45007
$ws.Cells.Item(9, 10).Value = 1160

$ws.Cells.Item(10, 4).Value = 45006

$ws.Cells.Item(11, 4).Value = 45070
$ws.Cells.Item(11, 10).Value = 800

$ws.Cells.Item(12, 4).Value = 45091
$ws.Cells.Item(12, 10).Value = 800

$ws.Cells.Item(13, 4).Value = 45028
$ws.Cells.Item(13, 10).Value = 1000

$ws.Cells.Item(14, 4).Value = 45084
$ws.Cells.Item(14, 10).Value = 900

$ws.Cells.Item(15, 4).Value = 44965
$ws.Cells.Item(15, 10).Value = 1120
$ws.Cells.Item(15, 11).Value = 2000
$ws.Cells.Item(15, 12).Value = 2500
$ws.Cells.Item(15, 13).Value = 2250
$ws.Cells.Item(15, 16).Value = 750

$ws.Cells.Item(16, 4).Value = 44911
$ws.Cells.Item(16, 10).Value = 700
$ws.Cells.Item(16, 11).Value = 1800
$ws.Cells.Item(16, 12).Value = 2000
$ws.Cells.Item(16, 13).Value = 1900
$ws.Cells.Item(16, 16).Value = 633

$ws.Cells.Item(17, 4).Value = 44971
$ws.Cells.Item(17, 10).Value = 1000

$ws.Cells.Item(18, 4).Value = 45035

$ws.Cells.Item(19, 4).Value = 44999
$ws.Cells.Item(19, 10).Value = 1100

$ws.Cells.Item(20, 4).Value = 44978
$ws.Cells.Item(20, 10).Value = 1000
$ws.Cells.Item(20, 11).Value = 1800
$ws.Cells.Item(20, 12).Value = 2000
$ws.Cells.Item(20, 13).Value = 1900
$ws.Cells.Item(20, 16).Value = 633

$ws.Cells.Item(21, 4).Value = 44985

$ws.Cells.Item(22, 4).Value = 45077
$ws.Cells.Item(22, 10).Value = 760

$ws.Cells.Item(23, 4).Value = 44910
$ws.Cells.Item(23, 10).Value = 1000
$ws.Cells.Item(23, 11).Value = 1800
$ws.Cells.Item(23, 13).Value = 1900
$ws.Cells.Item(23, 16).Value = 633

$ws.Cells.Item(24, 4).Value = 44953
$ws.Cells.Item(24, 10).Value = 1000

$ws.Cells.Item(25, 4).Value = 45062
$ws.Cells.Item(25, 10).Value = 1100
$ws.Cells.Item(25, 11).Value = 2000
$ws.Cells.Item(25, 12).Value = 2500
$ws.Cells.Item(25, 13).Value = 2250
$ws.Cells.Item(25, 16).Value = 750

$ws.Cells.Item(26, 4).Value = 44685
$ws.Cells.Item(26, 10).Value = 400
$ws.Cells.Item(26, 11).Value = 1500
$ws.Cells.Item(26, 12).Value = 2000
$ws.Cells.Item(26, 13).Value = 1750
$ws.Cells.Item(26, 16).Value = 583

$ws.Cells.Item(27, 4).Value = 44992
$ws.Cells.Item(27, 10).Value = 1040

$ws.Cells.Item(28, 4).Value = 44970
$ws.Cells.Item(28, 10).Value = 800

$ws.Cells.Item(29, 4).Value = 44883
$ws.Cells.Item(29, 10).Value = 500

$ws.Cells.Item(30, 4).Value = 44848
$ws.Cells.Item(30, 10).Value = 1000
$ws.Cells.Item(30, 11).Value = 1500
$ws.Cells.Item(30, 12).Value = 2000
$ws.Cells.Item(30, 13).Value = 1750
$ws.Cells.Item(30, 16).Value = 583

$ws.Cells.Item(31, 4).Value = 44964
